{"js": "// 1) Merge the three runs `{issue.` + `humanized_custom_if_need_additionally_contract_guarantee_issue_with_cost`\n//    + `}` in the same paragraph into a single run carrying the full placeholder text.\nconst mergedText =\n  \"{issue.humanized_custom_if_need_additionally_contract_guarantee_issue_with_cost}\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === mergedText) {\n    // Re-inserting the identical text with \"Replace\" collapses the paragraph's\n    // split runs into a single run.\n    para.getRange().insertText(mergedText, Word.InsertLocation.replace);\n    break;\n  }\n}\n\nawait context.sync();\n\n// 2) In the affiliates cell, the `|endfor}` placeholder should reference\n//    `issuer_affiliates_with_bank_liabilities` instead of `issuer_affiliates_all`\n//    (the `|for}` placeholder right before it stays untouched).\nconst searchResults = context.document.body.search(\n  \"{issue.issuer_affiliates_all|for}{issue.issuer_affiliates_all|endfor}\",\n  { matchCase: true }\n);\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"{issue.issuer_affiliates_all|for}{issue.issuer_affiliates_with_bank_liabilities|endfor}\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Merge the three runs `{issue.` + `humanized_custom_if_need_additionally_contract_guarantee_issue_with_cost`\n#    + `}` in the same paragraph into a single run carrying the full placeholder text.\n#    Re-running the same text through Find/Replace collapses the split runs into one.\n$mergedText = \"{issue.humanized_custom_if_need_additionally_contract_guarantee_issue_with_cost}\"\n\n$find1 = $d.Content.Find\n$find1.Text = $mergedText\n$find1.Replacement.Text = $mergedText\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# 2) In the affiliates cell, the `|endfor}` placeholder should reference\n#    `issuer_affiliates_with_bank_liabilities` instead of `issuer_affiliates_all`\n#    (the `|for}` placeholder right before it stays untouched).\n$find2 = $d.Content.Find\n$find2.Text = \"issuer_affiliates_all|endfor\"\n$find2.Replacement.Text = \"issuer_affiliates_with_bank_liabilities|endfor\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
